$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to reflect the new "through" date
$ws.Name = "Through 2022-03-19"

# Update the header label in I1 (shared string) to the new date
$ws.Range("I1").Value = "2022 (through 03-19)"

# Update the data values for the new day's data (2022-03-19)
$ws.Range("I4").Value = 84
$ws.Range("H6").Value = 109
$ws.Range("H14").Value = 1853
$ws.Range("I14").Value = 384
